$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 11 with new log entry: date, hours, task, source
# (copy B10's date formatting first so B11 keeps the existing date style
#  instead of Excel auto-creating a brand-new number format)
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B11").Value2 = 44605
$ws.Range("C11").Value2 = 1
$ws.Range("D11").Value = "Studio XGB "
$ws.Range("E11").Value = "Youtube "

# E12 picks up the underlined/centered style already used for D9
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the active selection to E12, matching where the user ended up
$ws.Range("E12").Select() | Out-Null
